# "Add password to LOGIN_USER"
# A new field row (password / VARCHAR2 / length 1024 / Nullable N) is inserted
# after the existing "role_id" row on the LOGIN_USER sheet. Because the sheet
# already has "spare" fully-styled blank rows below the data (through row 98),
# we reproduce the insert by shifting the existing field rows (6-10) down by
# one position each (writing directly into the cells rather than using a
# structural row-insert, so no new style records are created) and then
# populating the freed row 6 with the new "password" field.

$wb = $excel.ActiveWorkbook
$wsLogin    = $wb.Worksheets.Item("LOGIN_USER")
$wsContents = $wb.Worksheets.Item("Contents")

# --- Shift the field rows 10 -> 11, 9 -> 10, 8 -> 9, 7 -> 8, 6 -> 7 --------
# (iterate bottom-up so we never overwrite a row before reading it)
for ($r = 10; $r -ge 6; $r--) {
    $dst = $r + 1

    $bVal = $wsLogin.Cells.Item($r, 2).Value2
    $cVal = $wsLogin.Cells.Item($r, 3).Value2
    $dVal = $wsLogin.Cells.Item($r, 4).Value2
    $eVal = $wsLogin.Cells.Item($r, 5).Value2
    $fVal = $wsLogin.Cells.Item($r, 6).Value2

    $wsLogin.Cells.Item($dst, 1).Value = $r - 1          # NO. column keeps sequence 4..9
    $wsLogin.Cells.Item($dst, 2).Value = $bVal
    $wsLogin.Cells.Item($dst, 3).Value = $cVal

    if ($null -eq $dVal) {
        $wsLogin.Cells.Item($dst, 4).Clear()
    } else {
        $wsLogin.Cells.Item($dst, 4).Value = $dVal
    }

    $wsLogin.Cells.Item($dst, 5).Value = $eVal

    if ($null -eq $fVal) {
        $wsLogin.Cells.Item($dst, 6).Clear()
    } else {
        $wsLogin.Cells.Item($dst, 6).Value = $fVal
    }
}

# --- Populate the freed row 6 with the new "password" field ---------------
$wsLogin.Cells.Item(6, 1).Value = 4
$wsLogin.Cells.Item(6, 2).Value = "password"
$wsLogin.Cells.Item(6, 3).Value = "VARCHAR2"
$wsLogin.Cells.Item(6, 4).Value = 1024
$wsLogin.Cells.Item(6, 5).Value = "N"

# --- View-state bookkeeping, matching the authored selection/tab changes ---
# Contents was the selected tab before; LOGIN_USER becomes the selected tab.
[void]$wsContents.Activate()
[void]$wsContents.Range("B3").Select()

[void]$wsLogin.Activate()
[void]$wsLogin.Range("D6").Select()
